$p = $ppt.ActivePresentation

# 1) Slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS"): table on shape 2 switches
#    from the custom "Table_0" style to PowerPoint's built-in
#    "Medium Style 2 - Accent 1" table style.
$s5 = $p.Slides.Item(5)
$tbl = $s5.Shapes.Item(2).Table
$tbl.ApplyStyle("{4854F6EC-10AD-4A63-8430-1EB73821F3DE}")

# 2) The slide master's theme colour scheme is switched from the
#    "Integral" deck's Red Violet palette to the stock Office palette
#    (what used to live in the Notes Master's Office Theme).
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
